# Duplicate_Transactions.xlsx - "Now sending the repated message"
#
# 1) Fix B231 so the phone number is stored as a real number (it had been
#    left as text).
# 2) Append 7 new "Blocked" duplicate rows (232-238) recording that the
#    same book/message was resent to previously-blocked contacts.
#
# Phone numbers are numeric in this sheet except for the very last new
# row, which (matching the source export) keeps its phone number as text.
# A leading apostrophe forces Excel to store a numeric-looking value as
# text instead of re-interpreting it as a number; the same trick is used
# for the "Campaign_Date" column so dates like 2025-09-16 are kept as
# plain text rather than being converted into date serials.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) B231 was inline text "2065044242"; make it a genuine number ---
$ws.Range("B231").Value = 2065044242

# --- 2) Append the new duplicate rows ---
$newRows = @(
    @{ Row=232; A="David";            B=12814104622; BAsText=$false; C="PO Box 87301, Park Place, Houston, Texas";                 D="GG";  E="Po Box 87301, Park Place, Houston, Texas";                 L="2025-09-16 13:54:29" },
    @{ Row=233; A="Henry Chelegbor";  B=13024705411; BAsText=$false; C="6613 Guyer Street, Philadelphia, PA, Pennsylvania";         D="GG";  E="6613 Guyer Street, Philadelphia, Pa, Pennsylvania";         L="2025-09-16 13:54:31" },
    @{ Row=234; A="Dennis Vanmeter";  B=13049196111; BAsText=$false; C="1909 Harper Rd, Beckley, WV 25801";                         D="GG";  E="1909 Harper Rd, Beckley, Wv 25801";                         L="2025-09-16 13:54:33" },
    @{ Row=235; A="Madhukar Verma";   B=2065044242;  BAsText=$false; C="42729 Mayfair Park Ave Fremont Fremont 94538 California USA"; D="YBB"; E="English";                                                  L="2025-09-16 13:55:15" },
    @{ Row=236; A="David";            B=12814104622; BAsText=$false; C="PO Box 87301, Park Place, Houston, Texas";                 D="GG";  E="Po Box 87301, Park Place, Houston, Texas";                 L="2025-09-16 13:55:17" },
    @{ Row=237; A="Henry Chelegbor";  B=13024705411; BAsText=$false; C="6613 Guyer Street, Philadelphia, PA, Pennsylvania";         D="GG";  E="6613 Guyer Street, Philadelphia, Pa, Pennsylvania";         L="2025-09-16 13:55:20" },
    @{ Row=238; A="Dennis Vanmeter";  B="13049196111"; BAsText=$true; C="1909 Harper Rd, Beckley, WV 25801";                        D="GG";  E="1909 Harper Rd, Beckley, Wv 25801";                         L="2025-09-16 13:55:22" }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $ws.Range("A$row").Value = $r.A

    if ($r.BAsText) {
        # Force phone number to stay textual, e.g. "13049196111"
        $ws.Range("B$row").Value = "'" + $r.B
    } else {
        $ws.Range("B$row").Value = $r.B
    }

    $ws.Range("C$row").Value = $r.C
    $ws.Range("D$row").Value = $r.D
    $ws.Range("E$row").Value = $r.E
    # F, G, H, I, J are left blank for these rows.
    $ws.Range("K$row").Value = "Same book already sent"
    $ws.Range("L$row").Value = $r.L
    # Force the Campaign_Date to remain plain text "2025-09-16" rather
    # than becoming a date serial number.
    $ws.Range("M$row").Value = "'2025-09-16"
    $ws.Range("N$row").Value = "Blocked"
}
